$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every Price cell in column D holds a plain-text value in this sheet (e.g.
# "27.610.67", "1.00", "0.817") even though it looks numeric. Force each changed
# D-column cell to Text format before writing so Excel keeps it as text instead of
# reinterpreting/rounding it as a number (e.g. "1.00" -> 1, "65.30" -> 65.3).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '27.610.67'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '1.649.92'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '212.72'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D6").Value = '0.533'
$ws.Range("E6").Value = '  +4.86%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '23.66'
$ws.Range("E8").Value = '  -2.54%  '
$ws.Range("E9").Value = '  -1.42%  '
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("E11").Value = '  +1.63%  '
$ws.Range("D12").Value = '1.883.31'
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("D13").Value = '1.655.48'
$ws.Range("E13").Value = '  -0.40%  '
$ws.Range("E14").Value = '  +3.55%  '
$ws.Range("E15").Value = '  -2.22%  '
$ws.Range("D16").Value = '64.54'
$ws.Range("E16").Value = '  -2.25%  '
$ws.Range("D17").Value = '27.576.02'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '232.12'
$ws.Range("E18").Value = '  -3.63%  '
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("D20").Value = '7.59'
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D23").Value = '9.82'
$ws.Range("E23").Value = '  +4.25%  '
$ws.Range("E24").Value = '  -2.37%  '
$ws.Range("D25").Value = '148.61'
$ws.Range("E25").Value = '  +1.72%  '
$ws.Range("D26").Value = '7.02'
$ws.Range("E26").Value = '  -2.83%  '
$ws.Range("D27").Value = '0.113'
$ws.Range("E27").Value = '  +1.57%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("D29").Value = '15.66'
$ws.Range("E29").Value = '  -4.24%  '
$ws.Range("E30").Value = '  -2.59%  '
$ws.Range("D31").Value = '0.0488'
$ws.Range("E31").Value = '  -3.00%  '
$ws.Range("E32").Value = '  -0.84%  '
$ws.Range("D33").Value = '3.18'
$ws.Range("E33").Value = '  +1.97%  '
$ws.Range("D34").Value = '1.430.20'
$ws.Range("E34").Value = '  -2.24%  '
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("E37").Value = '  -0.64%  '
$ws.Range("E38").Value = '  -4.14%  '
$ws.Range("E39").Value = '  -3.42%  '
$ws.Range("D40").Value = '1.03'
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = '0.817'
$ws.Range("E42").Value = '  +3.15%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.55'
$ws.Range("E43").Value = '  +2.35%  '
$ws.Range("B44").Value = 'mCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D44").Value = '2.47'
$ws.Range("E44").Value = '  -3.03%  '
$ws.Range("E45").Value = '  +1.56%  '
$ws.Range("D46").Value = '65.30'
$ws.Range("E46").Value = '  -6.95%  '
$ws.Range("D47").Value = '1.792.79'
$ws.Range("E47").Value = '  -0.66%  '
$ws.Range("E48").Value = '  -1.19%  '
$ws.Range("D49").Value = '87.93'
$ws.Range("E49").Value = '  -0.98%  '
$ws.Range("D50").Value = '0.0₆0106'
$ws.Range("E50").Value = '  -2.58%  '
$ws.Range("E51").Value = '  -0.76%  '
